# resolved issue with date not excluding time
# The "Date" column (B) had accidentally been entered with a 12-day offset;
# correct the serial date values. Also remove the stray time-of-day /
# hours helper cells that leaked onto the grid (these carried the
# now-unused "h:mm:ss" time number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43599
$ws.Range("B3").Value = 43600
$ws.Range("B4").Value = 43601
$ws.Range("B5").Value = 43602
$ws.Range("B6").Value = 43603
$ws.Range("B7").Value = 43604
$ws.Range("B8").Value = 43605
$ws.Range("B11").Value = 43606
$ws.Range("B12").Value = 43607
$ws.Range("B13").Value = 43608
$ws.Range("B14").Value = 43609
$ws.Range("B15").Value = 43610
$ws.Range("B16").Value = 43611
$ws.Range("B17").Value = 43612

# Remove the leftover time/hours cells entirely (content + formatting),
# which also drops the now-unused time-of-day style from the workbook.
$ws.Range("G7").Clear()
$ws.Range("C11").Clear()
$ws.Range("D11").Clear()
$ws.Range("G11").Clear()
$ws.Range("C13").Clear()
$ws.Range("D13").Clear()
